$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp precision on the existing last row (row 7)
$ws.Cells.Item(7, 1).Value = 45865.25025472222

# Append a new row (row 8) with the latest scheduled-task reading
$ws.Cells.Item(8, 1).Value = 45865.3335573041
$ws.Cells.Item(8, 2).Value = 2025
$ws.Cells.Item(8, 3).Value = 30
$ws.Cells.Item(8, 4).Value = 13.83
$ws.Cells.Item(8, 5).Value = 89.66
$ws.Cells.Item(8, 6).Value = 59.77
$ws.Cells.Item(8, 7).Value = 1.98
$ws.Cells.Item(8, 8).Value = "N"
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = "08:00:19"

# Match the date-time format used by the rest of column A
$ws.Cells.Item(8, 1).NumberFormat = $ws.Cells.Item(7, 1).NumberFormat
